$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 135
$ws.Range("D5").Value = 22
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0.8598726114649682
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0.9246575342465754
